$wb = $excel.ActiveWorkbook

# ==== LP1912 ====
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:53:08"
$ws.Cells.Item(3, 1).Value = "Total filas: 282"
$ws.Cells.Item(39, 1).Value = "05:18:56"
$ws.Cells.Item(39, 3).Value = "215C_EL PATO"
$ws.Cells.Item(39, 4).Value = 88
$ws.Cells.Item(40, 1).Value = "06:43:40"
$ws.Cells.Item(40, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(40, 4).Value = 3
$ws.Cells.Item(47, 1).Value = "05:18:56"
$ws.Cells.Item(47, 3).Value = "15_ABASTO"
$ws.Cells.Item(47, 4).Value = 106
$ws.Cells.Item(48, 1).Value = "05:49:40"
$ws.Cells.Item(48, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(48, 4).Value = 75
$ws.Cells.Item(109, 1).Value = "07:59:28"
$ws.Cells.Item(109, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(109, 4).Value = 83
$ws.Cells.Item(110, 1).Value = "07:47:32"
$ws.Cells.Item(110, 3).Value = "17_ROMERO"
$ws.Cells.Item(110, 4).Value = 95
$ws.Cells.Item(117, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(118, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(145, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(146, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(189, 1).Value = "10:26:41"
$ws.Cells.Item(189, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(189, 4).Value = 100
$ws.Cells.Item(190, 1).Value = "10:56:30"
$ws.Cells.Item(190, 3).Value = "14_ABASTO"
$ws.Cells.Item(190, 4).Value = 70
$ws.Cells.Item(191, 1).Value = "12:01:50"
$ws.Cells.Item(191, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(191, 4).Value = 5
$ws.Cells.Item(192, 1).Value = "12:01:50"
$ws.Cells.Item(192, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(192, 4).Value = 5
$ws.Cells.Item(206, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(207, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(242, 1).Value = "13:19:56"
$ws.Cells.Item(242, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(242, 4).Value = 27
$ws.Cells.Item(243, 1).Value = "11:48:04"
$ws.Cells.Item(243, 3).Value = "17_ROMERO"
$ws.Cells.Item(243, 4).Value = 118
$ws.Cells.Item(251, 1).Value = "12:01:50"
$ws.Cells.Item(251, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(251, 4).Value = 115
$ws.Cells.Item(252, 1).Value = "12:37:14"
$ws.Cells.Item(252, 3).Value = "225_GOMEZ"
$ws.Cells.Item(252, 4).Value = 79
$ws.Cells.Item(254, 1).Value = "13:53:08"
$ws.Cells.Item(254, 2).Value = "13:57"
$ws.Cells.Item(254, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(254, 4).Value = 4
$ws.Cells.Item(255, 1).Value = "13:19:56"
$ws.Cells.Item(255, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(255, 4).Value = 45
$ws.Cells.Item(256, 1).Value = "12:37:14"
$ws.Cells.Item(256, 2).Value = "14:04"
$ws.Cells.Item(256, 3).Value = "17_ROMERO"
$ws.Cells.Item(256, 4).Value = 87
$ws.Cells.Item(257, 1).Value = "13:53:08"
$ws.Cells.Item(257, 2).Value = "14:05"
$ws.Cells.Item(257, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(257, 4).Value = 12
$ws.Cells.Item(258, 2).Value = "14:06"
$ws.Cells.Item(258, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(258, 4).Value = 47
$ws.Cells.Item(259, 1).Value = "12:55:01"
$ws.Cells.Item(259, 2).Value = "14:07"
$ws.Cells.Item(259, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(259, 4).Value = 72
$ws.Cells.Item(260, 1).Value = "13:53:08"
$ws.Cells.Item(260, 2).Value = "14:12"
$ws.Cells.Item(260, 3).Value = "15_ABASTO"
$ws.Cells.Item(260, 4).Value = 19
$ws.Cells.Item(261, 1).Value = "13:19:56"
$ws.Cells.Item(261, 2).Value = "14:16"
$ws.Cells.Item(261, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(261, 4).Value = 57
$ws.Cells.Item(262, 1).Value = "12:37:14"
$ws.Cells.Item(262, 2).Value = "14:17"
$ws.Cells.Item(262, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(262, 4).Value = 100
$ws.Cells.Item(263, 1).Value = "13:19:56"
$ws.Cells.Item(263, 2).Value = "14:19"
$ws.Cells.Item(263, 3).Value = "215C_EL PATO"
$ws.Cells.Item(263, 4).Value = 60
$ws.Cells.Item(264, 1).Value = "12:37:14"
$ws.Cells.Item(264, 2).Value = "14:20"
$ws.Cells.Item(264, 3).Value = "215C_EL PATO"
$ws.Cells.Item(264, 4).Value = 103
$ws.Cells.Item(265, 1).Value = "13:19:56"
$ws.Cells.Item(265, 2).Value = "14:20"
$ws.Cells.Item(265, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(265, 4).Value = 61
$ws.Cells.Item(266, 1).Value = "12:37:14"
$ws.Cells.Item(266, 2).Value = "14:21"
$ws.Cells.Item(266, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(266, 4).Value = 104
$ws.Cells.Item(267, 1).Value = "13:53:08"
$ws.Cells.Item(267, 2).Value = "14:28"
$ws.Cells.Item(267, 3).Value = "15_ABASTO"
$ws.Cells.Item(267, 4).Value = 35
$ws.Cells.Item(268, 1).Value = "13:53:08"
$ws.Cells.Item(268, 2).Value = "14:31"
$ws.Cells.Item(268, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(268, 4).Value = 38
$ws.Cells.Item(269, 1).Value = "13:53:08"
$ws.Cells.Item(269, 2).Value = "14:45"
$ws.Cells.Item(269, 3).Value = "14_ABASTO"
$ws.Cells.Item(269, 4).Value = 52
$ws.Cells.Item(270, 2).Value = "14:49"
$ws.Cells.Item(270, 3).Value = "14_ABASTO"
$ws.Cells.Item(270, 4).Value = 90
$ws.Cells.Item(271, 1).Value = "12:55:01"
$ws.Cells.Item(271, 2).Value = "14:50"
$ws.Cells.Item(271, 3).Value = "14_ABASTO"
$ws.Cells.Item(271, 4).Value = 115
$ws.Cells.Item(271, 5).Value = "LP1912"
$ws.Cells.Item(272, 1).Value = "13:19:56"
$ws.Cells.Item(272, 2).Value = "14:56"
$ws.Cells.Item(272, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(272, 4).Value = 97
$ws.Cells.Item(272, 5).Value = "LP1912"
$ws.Cells.Item(273, 1).Value = "13:53:08"
$ws.Cells.Item(273, 2).Value = "14:57"
$ws.Cells.Item(273, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(273, 4).Value = 64
$ws.Cells.Item(273, 5).Value = "LP1912"
$ws.Cells.Item(274, 1).Value = "13:19:56"
$ws.Cells.Item(274, 2).Value = "14:58"
$ws.Cells.Item(274, 3).Value = "215B_EL PATO"
$ws.Cells.Item(274, 4).Value = 99
$ws.Cells.Item(274, 5).Value = "LP1912"
$ws.Cells.Item(275, 1).Value = "13:19:56"
$ws.Cells.Item(275, 2).Value = "15:00"
$ws.Cells.Item(275, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(275, 4).Value = 101
$ws.Cells.Item(275, 5).Value = "LP1912"
$ws.Cells.Item(276, 1).Value = "13:19:56"
$ws.Cells.Item(276, 2).Value = "15:04"
$ws.Cells.Item(276, 3).Value = "10_OLMOS"
$ws.Cells.Item(276, 4).Value = 105
$ws.Cells.Item(276, 5).Value = "LP1912"
$ws.Cells.Item(277, 1).Value = "13:53:08"
$ws.Cells.Item(277, 2).Value = "15:05"
$ws.Cells.Item(277, 3).Value = "10_OLMOS"
$ws.Cells.Item(277, 4).Value = 72
$ws.Cells.Item(277, 5).Value = "LP1912"
$ws.Cells.Item(278, 1).Value = "13:53:08"
$ws.Cells.Item(278, 2).Value = "15:10"
$ws.Cells.Item(278, 3).Value = "17_ROMERO"
$ws.Cells.Item(278, 4).Value = 77
$ws.Cells.Item(278, 5).Value = "LP1912"
$ws.Cells.Item(279, 1).Value = "13:19:56"
$ws.Cells.Item(279, 2).Value = "15:13"
$ws.Cells.Item(279, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(279, 4).Value = 114
$ws.Cells.Item(279, 5).Value = "LP1912"
$ws.Cells.Item(280, 1).Value = "13:53:08"
$ws.Cells.Item(280, 2).Value = "15:14"
$ws.Cells.Item(280, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(280, 4).Value = 81
$ws.Cells.Item(280, 5).Value = "LP1912"
$ws.Cells.Item(281, 1).Value = "13:53:08"
$ws.Cells.Item(281, 2).Value = "15:28"
$ws.Cells.Item(281, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(281, 4).Value = 95
$ws.Cells.Item(281, 5).Value = "LP1912"
$ws.Cells.Item(282, 1).Value = "13:53:08"
$ws.Cells.Item(282, 2).Value = "15:32"
$ws.Cells.Item(282, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(282, 4).Value = 99
$ws.Cells.Item(282, 5).Value = "LP1912"
$ws.Cells.Item(283, 1).Value = "13:53:08"
$ws.Cells.Item(283, 2).Value = "15:35"
$ws.Cells.Item(283, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(283, 4).Value = 102
$ws.Cells.Item(283, 5).Value = "LP1912"
$ws.Cells.Item(284, 1).Value = "13:53:08"
$ws.Cells.Item(284, 2).Value = "15:37"
$ws.Cells.Item(284, 3).Value = "10_OLMOS"
$ws.Cells.Item(284, 4).Value = 104
$ws.Cells.Item(284, 5).Value = "LP1912"
$ws.Cells.Item(285, 1).Value = "13:53:08"
$ws.Cells.Item(285, 2).Value = "15:39"
$ws.Cells.Item(285, 3).Value = "215A_EL PATO"
$ws.Cells.Item(285, 4).Value = 106
$ws.Cells.Item(285, 5).Value = "LP1912"
$ws.Cells.Item(286, 1).Value = "13:53:08"
$ws.Cells.Item(286, 2).Value = "15:44"
$ws.Cells.Item(286, 3).Value = "14_ABASTO"
$ws.Cells.Item(286, 4).Value = 111
$ws.Cells.Item(286, 5).Value = "LP1912"
$ws.Cells.Item(287, 1).Value = "13:53:08"
$ws.Cells.Item(287, 2).Value = "15:47"
$ws.Cells.Item(287, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(287, 4).Value = 114
$ws.Cells.Item(287, 5).Value = "LP1912"

# ==== LP1912-215 ====
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:53:08"
$ws.Cells.Item(3, 1).Value = "Total filas: 32"
$ws.Cells.Item(37, 1).Value = "13:53:08"
$ws.Cells.Item(37, 2).Value = "15:39"
$ws.Cells.Item(37, 3).Value = "215A_EL PATO"
$ws.Cells.Item(37, 4).Value = 106
$ws.Cells.Item(37, 5).Value = "LP1912"

# ==== 6203-6173 ====
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:53:08"
$ws.Cells.Item(3, 1).Value = "Total filas: 43"
$ws.Cells.Item(48, 1).Value = "13:53:08"
$ws.Cells.Item(48, 2).Value = "15:34"
$ws.Cells.Item(48, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(48, 4).Value = 101
$ws.Cells.Item(48, 5).Value = "L6173"
